# commit: "new code and new analysis 2013_12_24_1:29"
#
# The six existing series-label rows (A3:A8) get an "LPS_" prefix, and the
# two labels they used to carry ("动_无动参" / "静_无动参") reappear as two
# brand-new rows appended at the bottom (A9:A10). The chuan*/单位 header
# cells and the numeric timing data are unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Relabel the existing series rows with the "LPS_" prefix.
$ws.Range("A3").Value = "LPS_动_无动参"
$ws.Range("A4").Value = "LPS_静_无动参"
$ws.Range("A5").Value = "LPS_动静_无动参"
$ws.Range("A6").Value = "LPS_动_有动参"
$ws.Range("A7").Value = "LPS_静_有动参"
$ws.Range("A8").Value = "LPS_动静_有动参"

# Append the two new rows with the original (un-prefixed) labels.
$ws.Range("A9").Value = "动_无动参"
$ws.Range("A10").Value = "静_无动参"

# Match the author's final cursor position.
$ws.Range("A10").Select() | Out-Null
